$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 100: Asking for a Friend / Beetle Glue (item 19906)
$ws.Range("H100").Value = 2000.7142
$ws.Range("I100").Value = 1251.25
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 1251.25
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -710.25
$ws.Range("N100").Value = -4082

# Row 135: For Tired Minds / Grade 1 Gemsap of Intelligence (item 44047)
$ws.Range("H135").Value = 1189.6
$ws.Range("I135").Value = 1623.6666
$ws.Range("J135").Value = 538.5
$ws.Range("K135").Value = 14612.9994
$ws.Range("L135").Value = 4846.5
$ws.Range("M135").Value = -12077.9994
$ws.Range("N135").Value = -9916.5

$ws = $wb.Worksheets.Item("ARM")
# Row 74: As the Bolt Flies / Titanium Nugget (item 44000)
$ws.Range("H74").Value = 957.4681
$ws.Range("I74").Value = 928.95123
$ws.Range("K74").Value = 928.95123
$ws.Range("M74").Value = -54.95123000000001

# Row 77: Heavy Metal Banned (L) / Titanium Nugget (item 44000)
$ws.Range("H77").Value = 957.4681
$ws.Range("I77").Value = 928.95123
$ws.Range("K77").Value = 4644.75615
$ws.Range("M77").Value = -276.7561500000002

$ws = $wb.Worksheets.Item("BSM")
# Row 62: Barring the Gates to Foundation / Mythrite Katzbalger (item 10586)
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# Row 65: Starting Young (L) / Mythrite Katzbalger (item 10586)
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# Row 99: Meddle in Metal / Oroshigane Ingot (item 19943)
$ws.Range("H99").Value = 58825750
$ws.Range("I99").Value = 83335680
$ws.Range("J99").Value = 1920
$ws.Range("K99").Value = 83335680
$ws.Range("L99").Value = 1920
$ws.Range("M99").Value = -83334182
$ws.Range("N99").Value = -4916

# Row 107: The Gold Experience / Deepgold Nugget (item 27706)
$ws.Range("H107").Value = 6289.7085
$ws.Range("I107").Value = 702.7895
$ws.Range("J107").Value = 27520
$ws.Range("K107").Value = 702.7895
$ws.Range("L107").Value = 27520
$ws.Range("M107").Value = 1217.2105
$ws.Range("N107").Value = -31360

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber (item 44023)
$ws.Range("H31").Value = 41825.19
$ws.Range("I31").Value = 3597.1875
$ws.Range("J31").Value = 102990
$ws.Range("K31").Value = 3597.1875
$ws.Range("L31").Value = 102990
$ws.Range("M31").Value = -3302.1875
$ws.Range("N31").Value = -103580

# Row 34: Armoires of the Rich and Famous / Walnut Lumber (item 44023)
$ws.Range("H34").Value = 41825.19
$ws.Range("I34").Value = 3597.1875
$ws.Range("J34").Value = 102990
$ws.Range("K34").Value = 3597.1875
$ws.Range("L34").Value = 102990
$ws.Range("M34").Value = -3395.1875
$ws.Range("N34").Value = -103394

# Row 44: Stay on Target / Yarzonshell Harpoon (item 1850)
$ws.Range("H44").Value = 21700
$ws.Range("J44").Value = 21700
$ws.Range("L44").Value = 21700
$ws.Range("N44").Value = -22584

# Row 58: You Do the Heavy Lifting / Mahogany Lumber (item 44021)
$ws.Range("H58").Value = 3186.6956
$ws.Range("I58").Value = 863.1579
$ws.Range("J58").Value = 14223.5
$ws.Range("K58").Value = 863.1579
$ws.Range("L58").Value = 14223.5
$ws.Range("M58").Value = -660.1579
$ws.Range("N58").Value = -14629.5

# Row 107: Built to Last / White Oak Lumber (item 27689)
$ws.Range("H107").Value = 394.14285
$ws.Range("I107").Value = 360.9091
$ws.Range("J107").Value = 430.7
$ws.Range("K107").Value = 360.9091
$ws.Range("L107").Value = 430.7
$ws.Range("M107").Value = 1559.0909
$ws.Range("N107").Value = -4270.7

# Row 132: Hull Lotta Damage / Ginseng Lumber (item 44019)
$ws.Range("H132").Value = 1417.5625
$ws.Range("I132").Value = 997.3913
$ws.Range("J132").Value = 2491.3333
$ws.Range("K132").Value = 2992.1739
$ws.Range("L132").Value = 7473.999899999999
$ws.Range("M132").Value = -462.1738999999998
$ws.Range("N132").Value = -12533.9999

# Row 134: Wood You Be Quiet / Ceiba Lumber (item 44020)
$ws.Range("H134").Value = 19231960
$ws.Range("I134").Value = 1138.35
$ws.Range("J134").Value = 83334696
$ws.Range("K134").Value = 3415.05
$ws.Range("L134").Value = 250004088
$ws.Range("M134").Value = -880.0499999999997
$ws.Range("N134").Value = -250009158

# Row 136: Turali Quality / Dark Mahogany Lumber (item 44021)
$ws.Range("H136").Value = 3186.6956
$ws.Range("I136").Value = 863.1579
$ws.Range("J136").Value = 14223.5
$ws.Range("K136").Value = 2589.4737
$ws.Range("L136").Value = 42670.5
$ws.Range("M136").Value = -39.47370000000001
$ws.Range("N136").Value = -47770.5

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap / Maple Syrup (item 43974)
$ws.Range("H5").Value = 1857.5385
$ws.Range("I5").Value = 225.38461
$ws.Range("J5").Value = 3489.6924
$ws.Range("K5").Value = 676.15383
$ws.Range("L5").Value = 10469.0772
$ws.Range("M5").Value = -564.15383
$ws.Range("N5").Value = -10693.0772

# Row 26: A Grape Idea / Grape Juice (item 4746)
$ws.Range("H26").Value = 270
$ws.Range("I26").Value = 100
$ws.Range("J26").Value = 298.33334
$ws.Range("K26").Value = 300
$ws.Range("L26").Value = 895.0000200000001
$ws.Range("M26").Value = -12
$ws.Range("N26").Value = -1471.00002

# Row 131: The Mountain Steeped / Tsai tou Vounou (item 36060)
$ws.Range("H131").Value = 20876344
$ws.Range("I131").Value = 62625350
$ws.Range("J131").Value = 1842.125
$ws.Range("K131").Value = 187876050
$ws.Range("L131").Value = 5526.375
$ws.Range("M131").Value = -187871010
$ws.Range("N131").Value = -15606.375

# Row 135: Not-so-secret Ingredient / Royal Maple Syrup (item 43974)
$ws.Range("H135").Value = 1857.5385
$ws.Range("I135").Value = 225.38461
$ws.Range("J135").Value = 3489.6924
$ws.Range("K135").Value = 2028.46149
$ws.Range("L135").Value = 31407.2316
$ws.Range("M135").Value = 506.5385099999999
$ws.Range("N135").Value = -36477.2316

# Row 140: Sweet, Sweet Bean Juice / Mesquite Juice (item 44097)
$ws.Range("H140").Value = 145072.86
$ws.Range("I140").Value = 152226.5
$ws.Range("K140").Value = 456679.5
$ws.Range("M140").Value = -451499.5

$ws = $wb.Worksheets.Item("GSM")
# Row 11: A Ringing Success / Copper Ring (item 4422)
$ws.Range("H11").Value = 10771777
$ws.Range("I11").Value = 10836091
$ws.Range("J11").Value = 10000000
$ws.Range("K11").Value = 10836091
$ws.Range("L11").Value = 10000000
$ws.Range("M11").Value = -10835952
$ws.Range("N11").Value = -10000278

# Row 21: Forever 21K / Brass Ring (item 4430)
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()

# Row 24: Bad Guys Eat Brass / Brass Ring of Crafting (item 4431)
$ws.Range("H24").Value = 106025000
$ws.Range("I24").Value = 151428580
$ws.Range("J24").Value = 83340.336
$ws.Range("K24").Value = 151428580
$ws.Range("L24").Value = 83340.336
$ws.Range("M24").Value = -151428407
$ws.Range("N24").Value = -83686.336

# Row 30: Dog Tags Are for Dogs / Brass Ring (item 4430)
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("N30").ClearContents()

# Row 70: Sky Is the Limit / Mythrite Ingot (item 14146)
$ws.Range("H70").Value = 4999.7
$ws.Range("I70").Value = 4999.5
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 4999.5
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -4729.5
$ws.Range("N70").Value = -5540

# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot (item 14146)
$ws.Range("H73").Value = 4999.7
$ws.Range("I73").Value = 4999.5
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 4999.5
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -4063.5
$ws.Range("N73").Value = -6872

# Row 132: On Board for Lar / Lar Ingot (item 44008)
$ws.Range("H132").Value = 2612.5557
$ws.Range("I132").Value = 2271.8333
$ws.Range("J132").Value = 4316.1665
$ws.Range("K132").Value = 6815.499899999999
$ws.Range("L132").Value = 12948.4995
$ws.Range("M132").Value = -4285.499899999999
$ws.Range("N132").Value = -18008.4995

$ws = $wb.Worksheets.Item("WVR")
# Row 113: A Tender Table / Pixie Floss (item 27752)
$ws.Range("H113").Value = 455.2381
$ws.Range("I113").Value = 452.42856
$ws.Range("J113").Value = 460.85715
$ws.Range("K113").Value = 1357.28568
$ws.Range("L113").Value = 1382.57145
$ws.Range("M113").Value = 812.71432
$ws.Range("N113").Value = -5722.571449999999
